$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '46.228.60'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.21%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.436.63'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +7.53%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.22%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '296.21'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.32%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '96.04'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.21%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.569'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.59%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.13%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.509'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.40%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.34'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.77%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0785'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.48%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.13'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.87%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.104'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.36%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.807.06'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +7.47%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.445.52'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +7.82%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.844'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +6.65%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.14'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.88%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '46.069.69'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.55%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.66'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.07%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0949'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.20%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.24'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +7.55%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.59'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.08%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '244.99'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.31%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.80'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.17%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.95'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +5.16%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '39.82'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.71%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.22'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.10%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.77'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.70%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.86'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +16.84%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '21.32'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.62%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.78'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.31%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '148.25'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.23%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.52'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.17%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0772'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.68%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.02'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +20.17%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.114'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.68%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.116'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.31%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '14.93'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.26%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.82'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.59%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0302'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.66%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.26'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.08%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.975.69'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +10.98%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.998'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.06%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '91.23'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.20%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.82'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.87%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '16.37'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +30.52%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.64'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +10.60%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '101.20'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +7.76%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.676.30'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +7.51%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.185'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.59%  '
